$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D20").Value = "[책 출간] 파이썬 생활 밀착형 프로젝트 (OpenAPI)"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/notice/561"

$ws.Range("D37").Value = "[Paper Review] LayoutLM: Pre-training of Text and Layout for Document Image Understanding"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1766&mod=document&pageid=1"

$ws.Range("D44").Value = "5G 관련주 분석(2) - HFR"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/82"
